# Fixed Casing and Spacing
# Re-assign the student names across the seating_plan grid (rows 2-4, columns A-J)
# to correct the previous row/column ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("seating_plan")

# Row 2
$ws.Range("A2").Value = "Lexi Green, poor"
$ws.Range("B2").Value = "Caitlin Boyd, poor"
$ws.Range("C2").Value = "Niko Morris, poor"
$ws.Range("D2").Value = "Stanley Hirst, poor"
$ws.Range("E2").Value = "James Calderon, poor"
$ws.Range("F2").Value = "Benedict Hobday, good"
$ws.Range("G2").Value = "Violet Hudson, poor"
$ws.Range("H2").Value = "Aarron Kelly, good"
$ws.Range("I2").Value = "James Shilton, poor"
$ws.Range("J2").Value = "Ava Lee, poor"

# Row 3
$ws.Range("B3").Value = "Benjamin Finn, good"
$ws.Range("C3").Value = "Katrina Petersone, good"
$ws.Range("D3").Value = "Brooke Layton, good"
$ws.Range("E3").Value = "Sophie Rayner, excellent"
$ws.Range("F3").Value = "Thomas Barrett, excellent"
$ws.Range("G3").Value = "Ruby Haigh, good"
$ws.Range("H3").Value = "William Hunt, good"
$ws.Range("I3").Value = "Nancy Enyoazu, good"
$ws.Range("J3").Value = "Madison Taylor, good"

# Row 4
$ws.Range("B4").Value = "Matthew Homan, excellent"
$ws.Range("C4").Value = "James Eilbeck, excellent"
$ws.Range("D4").Value = "Esther Sido, excellent"
$ws.Range("E4").Value = "Samuel Dixon, excellent"
$ws.Range("F4").Value = "Spencer Rowe, excellent"
$ws.Range("G4").Value = "Benjamin Hillary, excellent"
$ws.Range("H4").Value = "Alex Sentance, excellent"
